$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.356.59"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.251.11"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.632"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.42"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.81%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.623"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.16"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0951"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.21"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("E13").Value = "  -1.67%  "
$ws.Range("D14").Value = "2.590.65"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.855"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.08%  "
$ws.Range("D17").Value = "2.245.12"
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("D18").Value = "42.208.02"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "232.03"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +38.33%  "
$ws.Range("E25").Value = "  +0.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("E27").Value = "  -4.69%  "
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  +1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "169.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0825"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.53%  "
$ws.Range("E33").Value = "  +1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.28"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.97%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.125"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0316"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.93"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.60%  "
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "62.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.205"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "109.50"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("E46").Value = "  +1.32%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.998"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.20"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.69%  "
